# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (the quarter's fund-holding detail,
# mirroring the layout of the existing "2021-Q4" sheet) between "2021-Q4"
# and "总计", and adds a "2022-Q1" row at the top of "总计"'s summary table
# (shifting the existing "2021-Q4" row down one).

$wb = $excel.ActiveWorkbook

$sheetQ4 = $wb.Worksheets.Item("2021-Q4")
$sheetTotal = $wb.Worksheets.Item("总计")

# Recreate "总计" after a fresh "2022-Q1" sheet so sheetIds come out in
# document order (2021-Q4=1, 2022-Q1=2, 总计=3) rather than 总计 keeping its
# old id.
$sheetTotal.Delete()

$q1 = $wb.Worksheets.Add($null, $sheetQ4)
$q1.Name = "2022-Q1"

$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

# ---------------------------------------------------------------------
# "2022-Q1" detail sheet — same column layout as "2021-Q4":
# B=基金代码 C=基金名称 D=基金规模 E=股票总仓位 F=仓位占比 G=持有市值(亿元) H=仓位排名
# ---------------------------------------------------------------------

# Header row + index column (A2:A5), formatted like "2021-Q4".
$sheetQ4.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$sheetQ4.Range("A2").Copy()
$q1.Range("A2:A5").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# B2:G5 hold text that looks numeric ("014141", "0.52", "81.06", ...); force
# text format first so Excel doesn't coerce them into numbers (and drop
# leading zeros on the fund codes), then strip the number format back off so
# the cells end up with the sheet's default (no explicit) style, matching
# "2021-Q4"'s B:G data cells.
$q1.Range("B2:G5").NumberFormat = "@"

$q1.Range("B2").Value = "014141"
$q1.Range("C2").Value = "大成新能源混合A"
$q1.Range("D2").Value = "0.52"
$q1.Range("E2").Value = "81.06"
$q1.Range("F2").Value = "4.41"
$q1.Range("G2").Value = "0.0229"

$q1.Range("B3").Value = "009796"
$q1.Range("C3").Value = "大成汇享一年持有期混合A"
$q1.Range("D3").Value = "0.38"
$q1.Range("E3").Value = "22.99"
$q1.Range("F3").Value = "1.08"
$q1.Range("G3").Value = "0.0041"

$q1.Range("B4").Value = "014142"
$q1.Range("C4").Value = "大成新能源混合C"
$q1.Range("D4").Value = "0.02"
$q1.Range("E4").Value = "81.06"
$q1.Range("F4").Value = "4.41"
$q1.Range("G4").Value = "0.0009"

$q1.Range("B5").Value = "009797"
$q1.Range("C5").Value = "大成汇享一年持有期混合C"
$q1.Range("D5").Value = "0.04"
$q1.Range("E5").Value = "22.99"
$q1.Range("F5").Value = "1.08"
$q1.Range("G5").Value = "0.0004"

$q1.Range("B2:G5").ClearFormats()

# A2:A5 are the 0-based row index column; H2:H5 are plain numbers (rank).
$q1.Cells.Item(2, 1).Value = 0
$q1.Cells.Item(3, 1).Value = 1
$q1.Cells.Item(4, 1).Value = 2
$q1.Cells.Item(5, 1).Value = 3

$q1.Range("H2").Value = 8
$q1.Range("H3").Value = 9
$q1.Range("H4").Value = 8
$q1.Range("H5").Value = 9

# ---------------------------------------------------------------------
# "总计" summary sheet — unchanged header + new 2022-Q1 row + shifted
# 2021-Q4 row.
# B=日期 C=持有数量(只) D=持有市值(亿元)
# ---------------------------------------------------------------------

$sheetQ4.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$sheetQ4.Range("A2").Copy()
$total.Range("A2:A3").PasteSpecial(-4122)

$total.Range("B2:B3").NumberFormat = "@"
$total.Range("B2").Value = "2022-Q1"
$total.Range("B3").Value = "2021-Q4"
$total.Range("B2:B3").ClearFormats()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 3).Value = 4
$total.Cells.Item(2, 4).Value = 0.03

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 3).Value = 5
$total.Cells.Item(3, 4).Value = 1.16

$sheetQ4.Activate()
